# Update "Förändrad" (Changed) date column (C) for rows 2-7
# from serial date 45204 (2023-10-05) to serial date 45207 (2023-10-08)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newDate = Get-Date -Year 2023 -Month 10 -Day 8 -Hour 0 -Minute 0 -Second 0

foreach ($row in 2..7) {
    $ws.Cells.Item($row, 3).Value = $newDate
}
